$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 1507.1818
$ws.Range("J2").Value = 366.33334
$ws.Range("L2").Value = 366.33334
$ws.Range("N2").Value = -592.33334

# Row 76
$ws.Range("H76").Value = 6385
$ws.Range("I76").Value = 4424
$ws.Range("K76").Value = 4424
$ws.Range("M76").Value = -4109

# Row 79
$ws.Range("H79").Value = 6385
$ws.Range("I79").Value = 4424
$ws.Range("K79").Value = 4424
$ws.Range("M79").Value = -3332

# Row 116
$ws.Range("H116").Value = 5147.609
$ws.Range("I116").Value = 4512.625
$ws.Range("K116").Value = 4512.625
$ws.Range("M116").Value = -1070.625

# Row 138
$ws.Range("H138").Value = 2491.639
$ws.Range("I138").Value = 1167.85
$ws.Range("J138").Value = 3000.7886
$ws.Range("K138").Value = 3503.55
$ws.Range("L138").Value = 9002.3658
$ws.Range("M138").Value = 1636.45
$ws.Range("N138").Value = -19282.3658

$ws = $wb.Worksheets.Item("ARM")
# Row 6
$ws.Range("H6").Value = 10000
$ws.Range("I6").Value = 10000
$ws.Range("K6").Value = 10000
$ws.Range("M6").Value = -9827

# Row 25
$ws.Range("H25").Value = 1166.6666
$ws.Range("I25").Value = 500
$ws.Range("K25").Value = 500
$ws.Range("M25").Value = -98

# Row 61
$ws.Range("H61").Value = 50004836
$ws.Range("I61").Value = 38466736
$ws.Range("J61").Value = 125002500
$ws.Range("K61").Value = 38466736
$ws.Range("L61").Value = 125002500
$ws.Range("M61").Value = -38466524
$ws.Range("N61").Value = -125002924

# Row 74
$ws.Range("H74").Value = 11819178
$ws.Range("J74").Value = 1667703.6
$ws.Range("L74").Value = 1667703.6
$ws.Range("N74").Value = -1669451.6

# Row 77
$ws.Range("H77").Value = 11819178
$ws.Range("J77").Value = 1667703.6
$ws.Range("L77").Value = 8338518
$ws.Range("N77").Value = -8347254

# Row 103
$ws.Range("H103").Value = 51598
$ws.Range("J103").Value = 51598
$ws.Range("L103").Value = 51598
$ws.Range("N103").Value = -53942

# Row 122
$ws.Range("H122").Value = 3859.5881
$ws.Range("I122").Value = 3078.6
$ws.Range("K122").Value = 9235.799999999999
$ws.Range("M122").Value = -6785.799999999999

# Row 136
$ws.Range("H136").Value = 50004836
$ws.Range("I136").Value = 38466736
$ws.Range("J136").Value = 125002500
$ws.Range("K136").Value = 115400208
$ws.Range("L136").Value = 375007500
$ws.Range("M136").Value = -115397658
$ws.Range("N136").Value = -375012600

$ws = $wb.Worksheets.Item("BSM")
# Row 11
$ws.Range("H11").Value = 517.3333
$ws.Range("J11").Value = 500
$ws.Range("L11").Value = 500
$ws.Range("N11").Value = -780

# Row 86
$ws.Range("H86").Value = 2435.8
$ws.Range("I86").Value = 2641.4614
$ws.Range("K86").Value = 2641.4614
$ws.Range("M86").Value = -1518.4614

# Row 89
$ws.Range("H89").Value = 2435.8
$ws.Range("I89").Value = 2641.4614
$ws.Range("K89").Value = 13207.307
$ws.Range("M89").Value = -7591.307000000001

# Row 94
$ws.Range("H94").Value = 1300.4231
$ws.Range("I94").Value = 1111.8572
$ws.Range("K94").Value = 1111.8572
$ws.Range("M94").Value = -660.8571999999999

# Row 134
$ws.Range("H134").Value = 1668671.5
$ws.Range("I134").Value = 2403
$ws.Range("K134").Value = 7209
$ws.Range("M134").Value = -4674

$ws = $wb.Worksheets.Item("CRP")
# Row 4
$ws.Range("H4").Value = 4000000
$ws.Range("I4").Value = 4000000
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 4000000
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()
$ws.Range("M4").Value = -3999888

# Row 7
$ws.Range("H7").Value = 236.28572
$ws.Range("I7").Value = 211
$ws.Range("J7").Value = 299.5
$ws.Range("K7").Value = 211
$ws.Range("L7").Value = 299.5
$ws.Range("M7").Value = -98
$ws.Range("N7").Value = -525.5

# Row 23
$ws.Range("H23").Value = 2514.5
$ws.Range("I23").Value = 29
$ws.Range("J23").Value = 5000
$ws.Range("K23").Value = 29
$ws.Range("L23").Value = 5000
$ws.Range("M23").Value = 211
$ws.Range("N23").Value = -5480

# Row 27
$ws.Range("H27").Value = 2514.5
$ws.Range("I27").Value = 29
$ws.Range("J27").Value = 5000
$ws.Range("K27").Value = 29
$ws.Range("L27").Value = 5000
$ws.Range("M27").Value = 163
$ws.Range("N27").Value = -5384

# Row 105
$ws.Range("H105").Value = 1411.125
$ws.Range("I105").Value = 1276.8334
$ws.Range("K105").Value = 1276.8334
$ws.Range("M105").Value = 470.1666

$ws = $wb.Worksheets.Item("CUL")
# Row 8
$ws.Range("H8").Value = 350.5
$ws.Range("I8").Value = 350.5
$ws.Range("K8").Value = 1051.5
$ws.Range("M8").Value = -912.5

# Row 16
$ws.Range("H16").Value = 7501
$ws.Range("I16").Value = 3750.5
$ws.Range("K16").Value = 11251.5
$ws.Range("M16").Value = -11078.5

# Row 68
$ws.Range("H68").Value = 2687.7778
$ws.Range("I68").Value = 864.6667
$ws.Range("J68").Value = 3599.3333
$ws.Range("K68").Value = 2594.0001
$ws.Range("L68").Value = 10797.9999
$ws.Range("M68").Value = -1783.0001
$ws.Range("N68").Value = -12419.9999

# Row 71
$ws.Range("H71").Value = 2687.7778
$ws.Range("I71").Value = 864.6667
$ws.Range("J71").Value = 3599.3333
$ws.Range("K71").Value = 7782.0003
$ws.Range("L71").Value = 32393.9997
$ws.Range("M71").Value = -3726.0003
$ws.Range("N71").Value = -40505.9997

# Row 131
$ws.Range("H131").Value = 5154.1846
$ws.Range("I131").Value = 8931.666999999999
$ws.Range("J131").Value = 4770.0337
$ws.Range("K131").Value = 26795.001
$ws.Range("L131").Value = 14310.1011
$ws.Range("M131").Value = -21755.001
$ws.Range("N131").Value = -24390.1011

# Row 134
$ws.Range("H134").Value = 9527.429
$ws.Range("J134").Value = 10106.462
$ws.Range("L134").Value = 30319.386
$ws.Range("N134").Value = -40459.386

# Row 140
$ws.Range("H140").Value = 109546.79
$ws.Range("I140").Value = 113492.96
$ws.Range("K140").Value = 340478.88
$ws.Range("M140").Value = -335298.88

$ws = $wb.Worksheets.Item("GSM")
# Row 3
$ws.Range("H3").Value = 6263450.5
$ws.Range("I3").Value = 1000000
$ws.Range("J3").Value = 7316141
$ws.Range("K3").Value = 1000000
$ws.Range("L3").Value = 7316141
$ws.Range("M3").Value = -999884
$ws.Range("N3").Value = -7316373

# Row 7
$ws.Range("H7").Value = 2054001.2
$ws.Range("J7").Value = 85003
$ws.Range("L7").Value = 85003
$ws.Range("N7").Value = -85227

# Row 8
$ws.Range("H8").Value = 2054001.2
$ws.Range("J8").Value = 85003
$ws.Range("L8").Value = 85003
$ws.Range("N8").Value = -85281

# Row 14
$ws.Range("H14").Value = 10001500
$ws.Range("I14").Value = 3000
$ws.Range("J14").Value = 20000000
$ws.Range("K14").Value = 3000
$ws.Range("L14").Value = 20000000
$ws.Range("M14").Value = -2832
$ws.Range("N14").Value = -20000336

# Row 20
$ws.Range("H20").Value = 53227.332
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 53227.332
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 53227.332
$ws.Range("M20").ClearContents()
$ws.Range("N20").Value = -53717.332

# Row 22
$ws.Range("H22").Value = 2004
$ws.Range("I22").Value = 2004
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 2004
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -1475
$ws.Range("N22").ClearContents()

# Row 80
$ws.Range("H80").Value = 3577.6
$ws.Range("I80").Value = 3261.6365
$ws.Range("J80").Value = 4446.5
$ws.Range("K80").Value = 3261.6365
$ws.Range("L80").Value = 4446.5
$ws.Range("M80").Value = -2263.6365
$ws.Range("N80").Value = -6442.5

# Row 83
$ws.Range("H83").Value = 3577.6
$ws.Range("I83").Value = 3261.6365
$ws.Range("J83").Value = 4446.5
$ws.Range("K83").Value = 16308.1825
$ws.Range("L83").Value = 22232.5
$ws.Range("M83").Value = -11316.1825
$ws.Range("N83").Value = -32216.5

# Row 113
$ws.Range("H113").Value = 2760.1738
$ws.Range("I113").Value = 1564.9333
$ws.Range("J113").Value = 5001.25
$ws.Range("K113").Value = 1564.9333
$ws.Range("L113").Value = 5001.25
$ws.Range("M113").Value = 605.0667000000001
$ws.Range("N113").Value = -9341.25

# Row 126
$ws.Range("H126").Value = 4362.5713
$ws.Range("J126").Value = 4902.8
$ws.Range("L126").Value = 14708.4
$ws.Range("N126").Value = -19648.4

$ws = $wb.Worksheets.Item("LTW")
# Row 2
$ws.Range("H2").Value = 11148
$ws.Range("J2").Value = 11722
$ws.Range("L2").Value = 11722
$ws.Range("N2").Value = -11946

# Row 7
$ws.Range("H7").Value = 69334.56
$ws.Range("I7").Value = 5054.375
$ws.Range("K7").Value = 5054.375
$ws.Range("M7").Value = -4942.375

# Row 23
$ws.Range("H23").Value = 6
$ws.Range("I23").Value = 6
$ws.Range("K23").Value = 6
$ws.Range("M23").Value = 224

# Row 30
$ws.Range("H30").Value = 160
$ws.Range("I30").Value = 160
$ws.Range("K30").Value = 160
$ws.Range("M30").Value = -52

# Row 126
$ws.Range("H126").Value = 69334.56
$ws.Range("I126").Value = 5054.375
$ws.Range("K126").Value = 15163.125
$ws.Range("M126").Value = -12693.125

$ws = $wb.Worksheets.Item("WVR")
# Row 98
$ws.Range("H98").Value = 87222.5
$ws.Range("J98").Value = 87222.5
$ws.Range("L98").Value = 87222.5
$ws.Range("N98").Value = -93212.5

# Row 122
$ws.Range("H122").Value = 2924.7
$ws.Range("I122").Value = 2469.9412
$ws.Range("K122").Value = 7409.823600000001
$ws.Range("M122").Value = -4959.823600000001

# Row 123
$ws.Range("H123").Value = 76317.25
$ws.Range("J123").Value = 76317.25
$ws.Range("L123").Value = 76317.25
$ws.Range("N123").Value = -86117.25

# Row 125
$ws.Range("H125").Value = 22261.5
$ws.Range("J125").Value = 22261.5
$ws.Range("L125").Value = 22261.5
$ws.Range("N125").Value = -32101.5

# Row 132
$ws.Range("H132").Value = 21764.928
$ws.Range("J132").Value = 31624.25
$ws.Range("L132").Value = 94872.75
$ws.Range("N132").Value = -99932.75
